$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = [double]"1"
$ws.Range("F2").Value = [double]"0.3333333333333333"
$ws.Range("G2").Value = [double]"0.2189473333333334"
$ws.Range("H2").Value = [double]"0.656842"
$ws.Range("I2").Value = [double]"0.009402596261870986"
$ws.Range("J2").Value = [double]"0.009402596261870984"
$ws.Range("M2").Value = [double]"2.598166333333333"
$ws.Range("N2").Value = [double]"7.794499"
$ws.Range("O2").Value = [double]"0.3466013321552429"
$ws.Range("P2").Value = [double]"0.3466013321552429"
$ws.Range("Q2").Value = [double]"0.5688615902397779"
$ws.Range("R2").Value = [double]"5.119754312158"
$ws.Range("S2").Value = [double]"0.003258952390082391"
$ws.Range("T2").Value = [double]"0.00325895239008239"
$ws.Range("E3").Value = [double]"1"
$ws.Range("F3").Value = [double]"0.3333333333333333"
$ws.Range("G3").Value = [double]"0.2189473333333334"
$ws.Range("H3").Value = [double]"0.656842"
$ws.Range("I3").Value = [double]"0.009402596261870986"
$ws.Range("J3").Value = [double]"0.009402596261870984"
$ws.Range("M3").Value = [double]"4.333403333333333"
$ws.Range("O3").Value = [double]"0.5780859172985858"
$ws.Range("P3").Value = [double]"0.5780859172985858"
$ws.Range("Q3").Value = [double]"0.9487871040911111"
$ws.Range("R3").Value = [double]"8.539083936819999"
$ws.Range("S3").Value = [double]"0.005435508485031943"
$ws.Range("T3").Value = [double]"0.005435508485031942"
$ws.Range("E4").Value = [double]"1"
$ws.Range("F4").Value = [double]"0.3333333333333333"
$ws.Range("G4").Value = [double]"0.2189473333333334"
$ws.Range("H4").Value = [double]"0.656842"
$ws.Range("I4").Value = [double]"0.009402596261870986"
$ws.Range("J4").Value = [double]"0.009402596261870984"
$ws.Range("M4").Value = [double]"0.4692043333333333"
$ws.Range("N4").Value = [double]"1.407613"
$ws.Range("O4").Value = [double]"0.06259293136852516"
$ws.Range("P4").Value = [double]"0.06259293136852516"
$ws.Range("Q4").Value = [double]"0.1027310375717778"
$ws.Range("R4").Value = [double]"0.924579338146"
$ws.Range("S4").Value = [double]"0.0005885360625052418"
$ws.Range("T4").Value = [double]"0.0005885360625052417"
$ws.Range("E5").Value = [double]"1"
$ws.Range("F5").Value = [double]"0.3333333333333333"
$ws.Range("G5").Value = [double]"0.2189473333333334"
$ws.Range("H5").Value = [double]"0.656842"
$ws.Range("I5").Value = [double]"0.009402596261870986"
$ws.Range("J5").Value = [double]"0.009402596261870984"
$ws.Range("K5").Value = [double]"2"
$ws.Range("L5").Value = [double]"0.6666666666666666"
$ws.Range("M5").Value = [double]"0.09534933333333333"
$ws.Range("N5").Value = [double]"0.286048"
$ws.Range("O5").Value = [double]"0.01271981917764605"
$ws.Range("P5").Value = [double]"0.01271981917764604"
$ws.Range("Q5").Value = [double]"0.02087648226844445"
$ws.Range("R5").Value = [double]"0.187888340416"
$ws.Range("S5").Value = [double]"0.0001195993242514096"
$ws.Range("T5").Value = [double]"0.0001195993242514096"
$ws.Range("I6").Value = [double]"0.8622887582286424"
$ws.Range("J6").Value = [double]"0.8622887582286423"
$ws.Range("M6").Value = [double]"2.598166333333333"
$ws.Range("N6").Value = [double]"7.794499"
$ws.Range("O6").Value = [double]"0.3466013321552429"
$ws.Range("P6").Value = [double]"0.3466013321552429"
$ws.Range("Q6").Value = [double]"52.16888406035012"
$ws.Range("R6").Value = [double]"469.519956543151"
$ws.Range("S6").Value = [double]"0.2988704323045376"
$ws.Range("T6").Value = [double]"0.2988704323045376"
$ws.Range("I7").Value = [double]"0.8622887582286424"
$ws.Range("J7").Value = [double]"0.8622887582286423"
$ws.Range("M7").Value = [double]"4.333403333333333"
$ws.Range("O7").Value = [double]"0.5780859172985858"
$ws.Range("P7").Value = [double]"0.5780859172985858"
$ws.Range("Q7").Value = [double]"87.01090964925444"
$ws.Range("S7").Value = [double]"0.4984769877768632"
$ws.Range("T7").Value = [double]"0.4984769877768632"
$ws.Range("I8").Value = [double]"0.8622887582286424"
$ws.Range("J8").Value = [double]"0.8622887582286423"
$ws.Range("M8").Value = [double]"0.4692043333333333"
$ws.Range("N8").Value = [double]"1.407613"
$ws.Range("O8").Value = [double]"0.06259293136852516"
$ws.Range("P8").Value = [double]"0.06259293136852516"
$ws.Range("Q8").Value = [double]"9.421208393104113"
$ws.Range("R8").Value = [double]"84.790875537937"
$ws.Range("S8").Value = [double]"0.05397318106365619"
$ws.Range("T8").Value = [double]"0.05397318106365619"
$ws.Range("I9").Value = [double]"0.8622887582286424"
$ws.Range("J9").Value = [double]"0.8622887582286423"
$ws.Range("K9").Value = [double]"2"
$ws.Range("L9").Value = [double]"0.6666666666666666"
$ws.Range("M9").Value = [double]"0.09534933333333333"
$ws.Range("N9").Value = [double]"0.286048"
$ws.Range("O9").Value = [double]"0.01271981917764605"
$ws.Range("P9").Value = [double]"0.01271981917764604"
$ws.Range("Q9").Value = [double]"1.914530356305778"
$ws.Range("R9").Value = [double]"17.230773206752"
$ws.Range("S9").Value = [double]"0.01096815708358528"
$ws.Range("T9").Value = [double]"0.01096815708358528"
$ws.Range("G10").Value = [double]"2.823530666666667"
$ws.Range("H10").Value = [double]"8.470592"
$ws.Range("I10").Value = [double]"0.1212552739852724"
$ws.Range("J10").Value = [double]"0.1212552739852723"
$ws.Range("M10").Value = [double]"2.598166333333333"
$ws.Range("N10").Value = [double]"7.794499"
$ws.Range("O10").Value = [double]"0.3466013321552429"
$ws.Range("P10").Value = [double]"0.3466013321552429"
$ws.Range("Q10").Value = [double]"7.336002319267555"
$ws.Range("R10").Value = [double]"66.02402087340801"
$ws.Range("S10").Value = [double]"0.04202723949414437"
$ws.Range("T10").Value = [double]"0.04202723949414437"
$ws.Range("G11").Value = [double]"2.823530666666667"
$ws.Range("H11").Value = [double]"8.470592"
$ws.Range("I11").Value = [double]"0.1212552739852724"
$ws.Range("J11").Value = [double]"0.1212552739852723"
$ws.Range("M11").Value = [double]"4.333403333333333"
$ws.Range("O11").Value = [double]"0.5780859172985858"
$ws.Range("P11").Value = [double]"0.5780859172985858"
$ws.Range("Q11").Value = [double]"12.23549720270222"
$ws.Range("R11").Value = [double]"110.11947482432"
$ws.Range("S11").Value = [double]"0.07009596628906753"
$ws.Range("T11").Value = [double]"0.07009596628906752"
$ws.Range("G12").Value = [double]"2.823530666666667"
$ws.Range("H12").Value = [double]"8.470592"
$ws.Range("I12").Value = [double]"0.1212552739852724"
$ws.Range("J12").Value = [double]"0.1212552739852723"
$ws.Range("M12").Value = [double]"0.4692043333333333"
$ws.Range("N12").Value = [double]"1.407613"
$ws.Range("O12").Value = [double]"0.06259293136852516"
$ws.Range("P12").Value = [double]"0.06259293136852516"
$ws.Range("Q12").Value = [double]"1.324812824099556"
$ws.Range("R12").Value = [double]"11.923315416896"
$ws.Range("S12").Value = [double]"0.007589723042631867"
$ws.Range("T12").Value = [double]"0.007589723042631866"
$ws.Range("G13").Value = [double]"2.823530666666667"
$ws.Range("H13").Value = [double]"8.470592"
$ws.Range("I13").Value = [double]"0.1212552739852724"
$ws.Range("J13").Value = [double]"0.1212552739852723"
$ws.Range("K13").Value = [double]"2"
$ws.Range("L13").Value = [double]"0.6666666666666666"
$ws.Range("M13").Value = [double]"0.09534933333333333"
$ws.Range("N13").Value = [double]"0.286048"
$ws.Range("O13").Value = [double]"0.01271981917764605"
$ws.Range("P13").Value = [double]"0.01271981917764604"
$ws.Range("Q13").Value = [double]"0.2692217667128889"
$ws.Range("R13").Value = [double]"2.422995900416"
$ws.Range("S13").Value = [double]"0.001542345159428593"
$ws.Range("T13").Value = [double]"0.001542345159428593"
$ws.Range("G14").Value = [double]"0.1642436666666667"
$ws.Range("H14").Value = [double]"0.492731"
$ws.Range("I14").Value = [double]"0.007053371524214274"
$ws.Range("J14").Value = [double]"0.007053371524214274"
$ws.Range("M14").Value = [double]"2.598166333333333"
$ws.Range("N14").Value = [double]"7.794499"
$ws.Range("O14").Value = [double]"0.3466013321552429"
$ws.Range("P14").Value = [double]"0.3466013321552429"
$ws.Range("Q14").Value = [double]"0.4267323651965556"
$ws.Range("R14").Value = [double]"3.840591286769"
$ws.Range("S14").Value = [double]"0.002444707966478523"
$ws.Range("T14").Value = [double]"0.002444707966478523"
$ws.Range("G15").Value = [double]"0.1642436666666667"
$ws.Range("H15").Value = [double]"0.492731"
$ws.Range("I15").Value = [double]"0.007053371524214274"
$ws.Range("J15").Value = [double]"0.007053371524214274"
$ws.Range("M15").Value = [double]"4.333403333333333"
$ws.Range("O15").Value = [double]"0.5780859172985858"
$ws.Range("P15").Value = [double]"0.5780859172985858"
$ws.Range("Q15").Value = [double]"0.7117340526122222"
$ws.Range("R15").Value = [double]"6.40560647351"
$ws.Range("S15").Value = [double]"0.004077454747623133"
$ws.Range("T15").Value = [double]"0.004077454747623133"
$ws.Range("G16").Value = [double]"0.1642436666666667"
$ws.Range("H16").Value = [double]"0.492731"
$ws.Range("I16").Value = [double]"0.007053371524214274"
$ws.Range("J16").Value = [double]"0.007053371524214274"
$ws.Range("M16").Value = [double]"0.4692043333333333"
$ws.Range("N16").Value = [double]"1.407613"
$ws.Range("O16").Value = [double]"0.06259293136852516"
$ws.Range("P16").Value = [double]"0.06259293136852516"
$ws.Range("Q16").Value = [double]"0.07706384012255556"
$ws.Range("R16").Value = [double]"0.6935745611030001"
$ws.Range("S16").Value = [double]"0.0004414911997318537"
$ws.Range("T16").Value = [double]"0.0004414911997318537"
$ws.Range("G17").Value = [double]"0.1642436666666667"
$ws.Range("H17").Value = [double]"0.492731"
$ws.Range("I17").Value = [double]"0.007053371524214274"
$ws.Range("J17").Value = [double]"0.007053371524214274"
$ws.Range("K17").Value = [double]"2"
$ws.Range("L17").Value = [double]"0.6666666666666666"
$ws.Range("M17").Value = [double]"0.09534933333333333"
$ws.Range("N17").Value = [double]"0.286048"
$ws.Range("O17").Value = [double]"0.01271981917764605"
$ws.Range("P17").Value = [double]"0.01271981917764604"
$ws.Range("Q17").Value = [double]"0.01566052412088889"
$ws.Range("R17").Value = [double]"0.140944717088"
$ws.Range("S17").Value = [double]"8.971761038076323E-05"
$ws.Range("T17").Value = [double]"8.971761038076322E-05"
